$d = $word.ActiveDocument

# Merge the separate "<", "id>", "<id>.../id>" runs that spell out
# <id>p077r_4</id>, <id>p077v_1</id>, <id>p077v_2</id> and <id>p077v_3</id>
# into a single run/text-node each, as the XML was split across multiple
# <w:r> elements before and needs to become one run with one <w:t>.

$replacements = @(
    @{ old = "<id>p077r_4</id>"; new = "<id>p077r_4</id>" },
    @{ old = "<id>p077v_1</id>"; new = "<id>p077v_1</id>" },
    @{ old = "<id>p077v_2</id>"; new = "<id>p077v_2</id>" },
    @{ old = "<id>p077v_3</id>"; new = "<id>p077v_3</id>" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
